$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the symbol column (A2:A9) in an order that reproduces the
# shared-string table ordering seen in the target workbook (CRKN, AAPL,
# TSLA, INTC, GOOGL, META, LCID, NVDA appended in that order). A10 keeps
# referencing the pre-existing "mmm" string, unchanged.
$ws.Range("A9").Value = "CRKN"
$ws.Range("A8").Value = "AAPL"
$ws.Range("A7").Value = "TSLA"
$ws.Range("A6").Value = "INTC"
$ws.Range("A5").Value = "GOOGL"
$ws.Range("A4").Value = "META"
$ws.Range("A3").Value = "LCID"
$ws.Range("A2").Value = "NVDA"

# Update the date column with the corrected ("fixed") end dates.
$ws.Range("B2").Value = 45154
$ws.Range("B3").Value = 45154
$ws.Range("B4").Value = 45154
$ws.Range("B5").Value = 45155
$ws.Range("B6").Value = 45155
$ws.Range("B7").Value = 45155
$ws.Range("B8").Value = 45156
$ws.Range("B9").Value = 45156
$ws.Range("B10").Value = 45156

# Move the active selection to A15, matching the edited view state.
$ws.Range("A15").Select()
